$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(80, 8).Value = 785.56525  # H80: 1019.2941 -> 785.56525
$ws.Cells.Item(80, 9).Value = 264.66666  # I80: 346.57144 -> 264.66666
$ws.Cells.Item(80, 10).Value = 1353.8182  # J80: 1490.2 -> 1353.8182
$ws.Cells.Item(80, 11).Value = 793.9999799999999  # K80: 1039.71432 -> 793.9999799999999
$ws.Cells.Item(80, 12).Value = 4061.4546  # L80: 4470.6 -> 4061.4546
$ws.Cells.Item(80, 13).Value = 204.0000200000001  # M80: -41.71432000000004 -> 204.0000200000001
$ws.Cells.Item(80, 14).Value = -6057.4546  # N80: -6466.6 -> -6057.4546
$ws.Cells.Item(83, 8).Value = 785.56525  # H83: 1019.2941 -> 785.56525
$ws.Cells.Item(83, 9).Value = 264.66666  # I83: 346.57144 -> 264.66666
$ws.Cells.Item(83, 10).Value = 1353.8182  # J83: 1490.2 -> 1353.8182
$ws.Cells.Item(83, 11).Value = 2381.99994  # K83: 3119.14296 -> 2381.99994
$ws.Cells.Item(83, 12).Value = 12184.3638  # L83: 13411.8 -> 12184.3638
$ws.Cells.Item(83, 13).Value = 2610.00006  # M83: 1872.85704 -> 2610.00006
$ws.Cells.Item(83, 14).Value = -22168.3638  # N83: -23395.8 -> -22168.3638
$ws.Cells.Item(101, 8).Value = 474.3  # H101: 679.6 -> 474.3
$ws.Cells.Item(101, 9).Value = 471.625  # I101: 599.5 -> 471.625
$ws.Cells.Item(101, 10).Value = 485  # J101: 1000 -> 485
$ws.Cells.Item(101, 11).Value = 1414.875  # K101: 1798.5 -> 1414.875
$ws.Cells.Item(101, 12).Value = 1455  # L101: 3000 -> 1455
$ws.Cells.Item(101, 13).Value = 207.125  # M101: -176.5 -> 207.125
$ws.Cells.Item(101, 14).Value = -4699  # N101: -6244 -> -4699
$ws.Cells.Item(112, 8).Value = 1205.0385  # H112: 1218.125 -> 1205.0385
$ws.Cells.Item(112, 10).Value = 1230.4584  # J112: 1247.0454 -> 1230.4584
$ws.Cells.Item(112, 12).Value = 3691.3752  # L112: 3741.1362 -> 3691.3752
$ws.Cells.Item(112, 14).Value = -5907.3752  # N112: -5957.1362 -> -5907.3752
$ws.Cells.Item(116, 8).Value = 3137.6924  # H116: 2933.8572 -> 3137.6924
$ws.Cells.Item(116, 9).Value = 2310  # I116: 2370.3635 -> 2310
$ws.Cells.Item(116, 11).Value = 2310  # K116: 2370.3635 -> 2310
$ws.Cells.Item(116, 13).Value = 1132  # M116: 1071.6365 -> 1132
$ws.Cells.Item(134, 8).Value = 45695.555  # H134: 45686 -> 45695.555
$ws.Cells.Item(134, 10).Value = 45695.555  # J134: 45686 -> 45695.555
$ws.Cells.Item(134, 12).Value = 45695.555  # L134: 45686 -> 45695.555
$ws.Cells.Item(134, 14).Value = -55835.555  # N134: -55826 -> -55835.555
$ws.Cells.Item(137, 8).Value = 1203.1875  # H137: 893.5 -> 1203.1875
$ws.Cells.Item(137, 9).Value = 1332.5454  # I137: 800 -> 1332.5454
$ws.Cells.Item(137, 10).Value = 918.6  # J137: 949.6 -> 918.6
$ws.Cells.Item(137, 11).Value = 3997.6362  # K137: 2400 -> 3997.6362
$ws.Cells.Item(137, 12).Value = 2755.8  # L137: 2848.8 -> 2755.8
$ws.Cells.Item(137, 13).Value = -1447.6362  # M137: 150 -> -1447.6362
$ws.Cells.Item(137, 14).Value = -7855.8  # N137: -7948.8 -> -7855.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4221.27  # H32: 4503.34 -> 4221.27
$ws.Cells.Item(32, 9).Value = 4031.9082  # I32: 4354.165 -> 4031.9082
$ws.Cells.Item(32, 10).Value = 13500  # J32: 9326.666999999999 -> 13500
$ws.Cells.Item(32, 11).Value = 4031.9082  # K32: 4354.165 -> 4031.9082
$ws.Cells.Item(32, 12).Value = 13500  # L32: 9326.666999999999 -> 13500
$ws.Cells.Item(32, 13).Value = -3744.9082  # M32: -4067.165 -> -3744.9082
$ws.Cells.Item(32, 14).Value = -14074  # N32: -9900.666999999999 -> -14074

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(2, 8).Value = 0  # H2: 38856 -> 0
$ws.Cells.Item(2, 10).Value = 0  # J2: 38856 -> 0
$ws.Cells.Item(2, 12).Value = 0  # L2: 38856 -> 0
$ws.Cells.Item(2, 14).ClearContents()  # N2: delete (was -39082)
$ws.Cells.Item(86, 8).Value = 2407.4167  # H86: 3253.2856 -> 2407.4167
$ws.Cells.Item(86, 9).Value = 2053.5557  # I86: 3649.4285 -> 2053.5557
$ws.Cells.Item(86, 10).Value = 3469  # J86: 2857.1428 -> 3469
$ws.Cells.Item(86, 11).Value = 2053.5557  # K86: 3649.4285 -> 2053.5557
$ws.Cells.Item(86, 12).Value = 3469  # L86: 2857.1428 -> 3469
$ws.Cells.Item(86, 13).Value = -930.5556999999999  # M86: -2526.4285 -> -930.5556999999999
$ws.Cells.Item(86, 14).Value = -5715  # N86: -5103.1428 -> -5715
$ws.Cells.Item(89, 8).Value = 2407.4167  # H89: 3253.2856 -> 2407.4167
$ws.Cells.Item(89, 9).Value = 2053.5557  # I89: 3649.4285 -> 2053.5557
$ws.Cells.Item(89, 10).Value = 3469  # J89: 2857.1428 -> 3469
$ws.Cells.Item(89, 11).Value = 10267.7785  # K89: 18247.1425 -> 10267.7785
$ws.Cells.Item(89, 12).Value = 17345  # L89: 14285.714 -> 17345
$ws.Cells.Item(89, 13).Value = -4651.7785  # M89: -12631.1425 -> -4651.7785
$ws.Cells.Item(89, 14).Value = -28577  # N89: -25517.714 -> -28577

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 47881  # H28: 70643 -> 47881
$ws.Cells.Item(28, 10).Value = 47881  # J28: 70643 -> 47881
$ws.Cells.Item(28, 12).Value = 47881  # L28: 70643 -> 47881
$ws.Cells.Item(28, 14).Value = -48371  # N28: -71133 -> -48371
$ws.Cells.Item(31, 8).Value = 1785.31  # H31: 1954.0505 -> 1785.31
$ws.Cells.Item(31, 9).Value = 1153.6984  # I31: 1233.3442 -> 1153.6984
$ws.Cells.Item(31, 10).Value = 2860.7568  # J31: 3110.9736 -> 2860.7568
$ws.Cells.Item(31, 11).Value = 1153.6984  # K31: 1233.3442 -> 1153.6984
$ws.Cells.Item(31, 12).Value = 2860.7568  # L31: 3110.9736 -> 2860.7568
$ws.Cells.Item(31, 13).Value = -858.6984  # M31: -938.3442 -> -858.6984
$ws.Cells.Item(31, 14).Value = -3450.7568  # N31: -3700.9736 -> -3450.7568
$ws.Cells.Item(34, 8).Value = 1785.31  # H34: 1954.0505 -> 1785.31
$ws.Cells.Item(34, 9).Value = 1153.6984  # I34: 1233.3442 -> 1153.6984
$ws.Cells.Item(34, 10).Value = 2860.7568  # J34: 3110.9736 -> 2860.7568
$ws.Cells.Item(34, 11).Value = 1153.6984  # K34: 1233.3442 -> 1153.6984
$ws.Cells.Item(34, 12).Value = 2860.7568  # L34: 3110.9736 -> 2860.7568
$ws.Cells.Item(34, 13).Value = -951.6984  # M34: -1031.3442 -> -951.6984
$ws.Cells.Item(34, 14).Value = -3264.7568  # N34: -3514.9736 -> -3264.7568

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(34, 8).Value = 697.2273  # H34: 664.125 -> 697.2273
$ws.Cells.Item(34, 10).Value = 736.95  # J34: 697.2273 -> 736.95
$ws.Cells.Item(34, 12).Value = 2210.85  # L34: 2091.6819 -> 2210.85
$ws.Cells.Item(34, 14).Value = -2378.85  # N34: -2259.6819 -> -2378.85
$ws.Cells.Item(107, 8).Value = 545087.0600000001  # H107: 486726.78 -> 545087.0600000001
$ws.Cells.Item(107, 9).Value = 788.6  # I107: 612.7778 -> 788.6
$ws.Cells.Item(107, 10).Value = 681161.6  # J107: 716991.3 -> 681161.6
$ws.Cells.Item(107, 11).Value = 2365.8  # K107: 1838.3334 -> 2365.8
$ws.Cells.Item(107, 12).Value = 2043484.8  # L107: 2150973.9 -> 2043484.8
$ws.Cells.Item(107, 13).Value = -445.8000000000002  # M107: 81.66660000000002 -> -445.8000000000002
$ws.Cells.Item(107, 14).Value = -2047324.8  # N107: -2154813.9 -> -2047324.8
$ws.Cells.Item(122, 8).Value = 958.5  # H122: 848.5135 -> 958.5
$ws.Cells.Item(122, 9).Value = 530.36365  # I122: 460.39285 -> 530.36365
$ws.Cells.Item(122, 10).Value = 1900.4  # J122: 2056 -> 1900.4
$ws.Cells.Item(122, 11).Value = 4773.27285  # K122: 4143.53565 -> 4773.27285
$ws.Cells.Item(122, 12).Value = 17103.6  # L122: 18504 -> 17103.6
$ws.Cells.Item(122, 13).Value = -2323.27285  # M122: -1693.53565 -> -2323.27285
$ws.Cells.Item(122, 14).Value = -22003.6  # N122: -23404 -> -22003.6
$ws.Cells.Item(129, 8).Value = 1660.6333  # H129: 1470.5294 -> 1660.6333
$ws.Cells.Item(129, 9).Value = 892  # I129: 673.26666 -> 892
$ws.Cells.Item(129, 10).Value = 2044.95  # J129: 2099.9473 -> 2044.95
$ws.Cells.Item(129, 11).Value = 2676  # K129: 2019.79998 -> 2676
$ws.Cells.Item(129, 12).Value = 6134.85  # L129: 6299.841899999999 -> 6134.85
$ws.Cells.Item(129, 13).Value = 2324  # M129: 2980.20002 -> 2324
$ws.Cells.Item(129, 14).Value = -16134.85  # N129: -16299.8419 -> -16134.85
$ws.Cells.Item(131, 8).Value = 896.33  # H131: 853.22 -> 896.33
$ws.Cells.Item(131, 9).Value = 457.5  # I131: 288 -> 457.5
$ws.Cells.Item(131, 10).Value = 914.61456  # J131: 882.96844 -> 914.61456
$ws.Cells.Item(131, 11).Value = 1372.5  # K131: 864 -> 1372.5
$ws.Cells.Item(131, 12).Value = 2743.84368  # L131: 2648.90532 -> 2743.84368
$ws.Cells.Item(131, 13).Value = 3667.5  # M131: 4176 -> 3667.5
$ws.Cells.Item(131, 14).Value = -12823.84368  # N131: -12728.90532 -> -12823.84368
$ws.Cells.Item(132, 8).Value = 1980.0435  # H132: 2178.5789 -> 1980.0435
$ws.Cells.Item(132, 9).Value = 1599.5714  # I132: 2074.75 -> 1599.5714
$ws.Cells.Item(132, 10).Value = 2146.5  # J132: 2206.2666 -> 2146.5
$ws.Cells.Item(132, 11).Value = 14396.1426  # K132: 18672.75 -> 14396.1426
$ws.Cells.Item(132, 12).Value = 19318.5  # L132: 19856.3994 -> 19318.5
$ws.Cells.Item(132, 13).Value = -11866.1426  # M132: -16142.75 -> -11866.1426
$ws.Cells.Item(132, 14).Value = -24378.5  # N132: -24916.3994 -> -24378.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 76485.71000000001  # H40: 105800 -> 76485.71000000001
$ws.Cells.Item(40, 9).Value = 130100  # I40: 173000 -> 130100
$ws.Cells.Item(40, 11).Value = 130100  # K40: 173000 -> 130100
$ws.Cells.Item(40, 13).Value = -129964  # M40: -172864 -> -129964
$ws.Cells.Item(46, 8).Value = 949.7778  # H46: 1035.7222 -> 949.7778
$ws.Cells.Item(46, 9).Value = 880.4286  # I46: 876.2 -> 880.4286
$ws.Cells.Item(46, 10).Value = 1192.5  # J46: 1833.3334 -> 1192.5
$ws.Cells.Item(46, 11).Value = 880.4286  # K46: 876.2 -> 880.4286
$ws.Cells.Item(46, 12).Value = 1192.5  # L46: 1833.3334 -> 1192.5
$ws.Cells.Item(46, 13).Value = -692.4286  # M46: -688.2 -> -692.4286
$ws.Cells.Item(46, 14).Value = -1568.5  # N46: -2209.3334 -> -1568.5
$ws.Cells.Item(81, 8).Value = 30000  # H81: 0 -> 30000
$ws.Cells.Item(81, 10).Value = 30000  # J81: 0 -> 30000
$ws.Cells.Item(81, 12).Value = 30000  # L81: 0 -> 30000
$ws.Cells.Item(81, 14).Value = -31996  # N81: None -> -31996
$ws.Cells.Item(84, 8).Value = 30000  # H84: 0 -> 30000
$ws.Cells.Item(84, 10).Value = 30000  # J84: 0 -> 30000
$ws.Cells.Item(84, 12).Value = 90000  # L84: 0 -> 90000
$ws.Cells.Item(84, 14).Value = -99984  # N84: None -> -99984
$ws.Cells.Item(93, 8).Value = 854.1818  # H93: 881.5454999999999 -> 854.1818
$ws.Cells.Item(93, 9).Value = 874  # I93: 911.625 -> 874
$ws.Cells.Item(93, 11).Value = 874  # K93: 911.625 -> 874
$ws.Cells.Item(93, 13).Value = 374  # M93: 336.375 -> 374
$ws.Cells.Item(132, 8).Value = 5741.0986  # H132: 5628.473 -> 5741.0986
$ws.Cells.Item(132, 9).Value = 5386.316  # I132: 5350.6724 -> 5386.316
$ws.Cells.Item(132, 10).Value = 7185.5713  # J132: 6635.5 -> 7185.5713
$ws.Cells.Item(132, 11).Value = 16158.948  # K132: 16052.0172 -> 16158.948
$ws.Cells.Item(132, 12).Value = 21556.7139  # L132: 19906.5 -> 21556.7139
$ws.Cells.Item(132, 13).Value = -13628.948  # M132: -13522.0172 -> -13628.948
$ws.Cells.Item(132, 14).Value = -26616.7139  # N132: -24966.5 -> -26616.7139
$ws.Cells.Item(136, 8).Value = 13336174  # H136: 15154429 -> 13336174
$ws.Cells.Item(136, 9).Value = 2737.1052  # I136: 2838.6667 -> 2737.1052
$ws.Cells.Item(136, 10).Value = 55558724  # J136: 83336584 -> 55558724
$ws.Cells.Item(136, 11).Value = 8211.3156  # K136: 8516.000100000001 -> 8211.3156
$ws.Cells.Item(136, 12).Value = 166676172  # L136: 250009752 -> 166676172
$ws.Cells.Item(136, 13).Value = -5661.3156  # M136: -5966.000100000001 -> -5661.3156
$ws.Cells.Item(136, 14).Value = -166681272  # N136: -250014852 -> -166681272

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 2195.4075  # H81: 2281.0625 -> 2195.4075
$ws.Cells.Item(81, 9).Value = 2100.5386  # I81: 1599.5 -> 2100.5386
$ws.Cells.Item(81, 10).Value = 2283.5  # J81: 2690 -> 2283.5
$ws.Cells.Item(81, 11).Value = 4201.0772  # K81: 3199 -> 4201.0772
$ws.Cells.Item(81, 12).Value = 4567  # L81: 5380 -> 4567
$ws.Cells.Item(81, 13).Value = -3140.0772  # M81: -2138 -> -3140.0772
$ws.Cells.Item(81, 14).Value = -6689  # N81: -7502 -> -6689
$ws.Cells.Item(84, 8).Value = 2195.4075  # H84: 2281.0625 -> 2195.4075
$ws.Cells.Item(84, 9).Value = 2100.5386  # I84: 1599.5 -> 2100.5386
$ws.Cells.Item(84, 10).Value = 2283.5  # J84: 2690 -> 2283.5
$ws.Cells.Item(84, 11).Value = 21005.386  # K84: 15995 -> 21005.386
$ws.Cells.Item(84, 12).Value = 22835  # L84: 26900 -> 22835
$ws.Cells.Item(84, 13).Value = -15701.386  # M84: -10691 -> -15701.386
$ws.Cells.Item(84, 14).Value = -33443  # N84: -37508 -> -33443
